$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARTW")

# Insert a new column before column D, shifting existing D:K to E:L
$ws.Range("D1").EntireColumn.Insert()

Write-Host "Inserted column"
